# Add "docs" marker column (K) to the Tridev PCB BOM sheet.
# Every data row (2-33) gets a "y" in column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 33; $row++) {
    $ws.Range("K" + $row).Value = "y"
}

# Match the author's final cursor/scroll position.
$ws.Range("D18").Select()
